$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-33: update date serial value 45224 -> 45233
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
